$wb = $excel.ActiveWorkbook

# --- Metadata sheet: fix URL casing (corecl -> CoreCL) and bump the Date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "https://hl7chile.cl/fhir/ig/CoreCL/StructureDefinition/NotasEncuentro"
$meta.Range("B8").Value = "2022-12-12T20:08:16-03:00"

# --- Elements sheet: column widths shrank (re-fit after the cleanup pass) ---
$ws = $wb.Worksheets.Item("Elements")
$ws.Columns.Item(1).ColumnWidth = 18.166666666666668
$ws.Columns.Item(2).ColumnWidth = 10.333333333333334
$ws.Columns.Item(3).ColumnWidth = 6.833333333333333
$ws.Columns.Item(4).ColumnWidth = 5.0
$ws.Columns.Item(5).ColumnWidth = 3.8333333333333335
$ws.Columns.Item(6).ColumnWidth = 4.166666666666667
$ws.Columns.Item(7).ColumnWidth = 13.833333333333334
$ws.Columns.Item(8).ColumnWidth = 11.166666666666666
$ws.Columns.Item(9).ColumnWidth = 11.833333333333334
$ws.Columns.Item(11).ColumnWidth = 40.666666666666664
$ws.Columns.Item(15).ColumnWidth = 12.666666666666666
$ws.Columns.Item(20).ColumnWidth = 14.833333333333334
$ws.Columns.Item(21).ColumnWidth = 15.333333333333334
$ws.Columns.Item(22).ColumnWidth = 16.166666666666668
$ws.Columns.Item(23).ColumnWidth = 15.5
$ws.Columns.Item(24).ColumnWidth = 18.0
$ws.Columns.Item(25).ColumnWidth = 16.166666666666668
$ws.Columns.Item(26).ColumnWidth = 4.833333333333333
$ws.Columns.Item(27).ColumnWidth = 18.833333333333332
$ws.Columns.Item(28).ColumnWidth = 39.166666666666664
$ws.Columns.Item(29).ColumnWidth = 14.166666666666666
$ws.Columns.Item(30).ColumnWidth = 11.5
$ws.Columns.Item(31).ColumnWidth = 16.833333333333332
$ws.Columns.Item(32).ColumnWidth = 8.666666666666666
$ws.Columns.Item(33).ColumnWidth = 9.0
$ws.Columns.Item(34).ColumnWidth = 11.333333333333334
$ws.Columns.Item(36).ColumnWidth = 21.833333333333332
